# "Filled out excel document" -- populate the Dashboard sheet of the
# BoardGame tracking workbook with the team's actual data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dashboard")

# --- Header block (team name / game name, and the two side summary boxes) ---
$ws.Range("C4").Value = "Team Name: T3, aka"
$ws.Range("D4").Value = "Spaghetti Monster Acolytes"
$ws.Range("I4").Value = "Member Contribution"
$ws.Range("L4").Value = "Component Effort"

$ws.Range("D5").Value = "Bored Game"

$ws.Range("I6").Value = "Team Member Name"
$ws.Range("J6").Value = "Total time Spent (hours)"
$ws.Range("L6").Value = "Component Name"
$ws.Range("M6").Value = "Total time spent on component (hours)"

# --- Understanding the problem (row 8) ---
$ws.Range("D8").Value = "4 hours"
$ws.Range("E8").Value = "10 times"
$ws.Range("F8").Value = "see gitHub"
$ws.Range("G8").Value = 100
$ws.Range("I8").Value = "Adam Stammer"
$ws.Range("J8").Value = 24.25
$ws.Range("L8").Value = "Network Handler"
$ws.Range("M8").Value = 7.25

# --- Creating the model (row 10) ---
$ws.Range("D10").Value = "4 hours sketching diagrams"
$ws.Range("E10").Value = "15 times"
$ws.Range("F10").Value = "see gitHub"
$ws.Range("G10").Value = 100
$ws.Range("I10").Value = "Henry Weber"
$ws.Range("J10").Value = 8.75
$ws.Range("L10").Value = "Game"
$ws.Range("M10").Value = 21.25

# --- Writing code (row 12) ---
$ws.Range("D12").Value = "32.0 hours"
$ws.Range("E12").Value = "30 times"
$ws.Range("F12").Value = "see gitHub"
$ws.Range("G12").Value = "95 (lacking sanitation and error recovery)"
$ws.Range("I12").Value = "Stephanie Smith"
$ws.Range("J12").Value = 11.25
$ws.Range("L12").Value = "Gui"
$ws.Range("M12").Value = 11.25

# --- Testing the code (row 14) ---
$ws.Range("D14").Value = "12 hours"
$ws.Range("E14").Value = "100+ times"
$ws.Range("F14").Value = "see gitHub"
$ws.Range("G14").Value = "90 (could always use more testing)"
$ws.Range("I14").Value = "Akin Tema-Lopez"
$ws.Range("J14").Value = 5.75
$ws.Range("L14").Value = "Match"
$ws.Range("M14").Value = 13.25

# --- Feedback from testing and completion (row 16) ---
$ws.Range("D16").Value = "1 hours talking about it after testing"
$ws.Range("E16").Value = "20+ times"
$ws.Range("F16").Value = "camaraderie"
$ws.Range("G16").Value = 100
$ws.Range("M16").Value = "total 53 hours"

# New row 17: the team's running total, directly below the table.
$ws.Range("D17").Value = "total 53 hours"

# --- Row heights / column widths to fit the new, taller wrapped text ---
$ws.Columns.Item(3).ColumnWidth = 15
$ws.Rows.Item(4).RowHeight = 32
$ws.Rows.Item(6).RowHeight = 64
$ws.Rows.Item(10).RowHeight = 48
$ws.Rows.Item(14).RowHeight = 48
$ws.Rows.Item(17).RowHeight = 16

# Final cursor position left where the author's last edit landed.
$ws.Range("L16").Select()
